# Refactor synthetic array /3: swap color-square emojis for color-book emojis
# and rename the "noir" (black) status label to "bleu" (blue).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map old emoji/label -> new emoji/label
$emojiMap = @{
    "🟥" = "📕"  # red square   -> red book
    "⬛" = "📘"  # black square -> blue book
    "🟧" = "📙"  # orange square -> orange book
    "🟩" = "📗"  # green square -> green book
}

$labelMap = @{
    "noir" = "bleu"
}

$used = $ws.UsedRange
$rows = $used.Rows.Count

for ($r = 1; $r -le $rows; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $valA = $cellA.Value2
    if ($emojiMap.ContainsKey($valA)) {
        $cellA.Value2 = $emojiMap[$valA]
    }

    $cellB = $ws.Cells.Item($r, 2)
    $valB = $cellB.Value2
    if ($labelMap.ContainsKey($valB)) {
        $cellB.Value2 = $labelMap[$valB]
    }
}
